$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "1.003", "22.171.79").
# Force it to stay text so Excel does not coerce it into a Number/Date.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.171.79"
$ws.Range("E2").Value = "  -1.24%  "

$ws.Range("D3").Value = "1.554.51"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").Value = "288.03"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("D7").Value = "0.3823"
$ws.Range("E7").Value = "  +2.76%  "

$ws.Range("D8").Value = "0.3308"
$ws.Range("E8").Value = "  -0.29%  "

$ws.Range("D9").Value = "43.77"
$ws.Range("E9").Value = "  -9.27%  "

$ws.Range("D10").Value = "1.135"
$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("D11").Value = "0.07368"
$ws.Range("E11").Value = "  -1.35%  "

$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").Value = "20.14"
$ws.Range("E13").Value = "  -2.87%  "

$ws.Range("D14").Value = "5.816"
$ws.Range("E14").Value = "  -2.35%  "

$ws.Range("D15").Value = "1.580.34"
$ws.Range("E15").Value = "  -0.02%  "

$ws.Range("D16").Value = "6.729"
$ws.Range("E16").Value = "  -2.50%  "

$ws.Range("D17").Value = "0.00001072"
$ws.Range("E17").Value = "  -3.71%  "

$ws.Range("D18").Value = "0.06657"
$ws.Range("E18").Value = "  -1.61%  "

$ws.Range("D19").Value = "85.82"
$ws.Range("E19").Value = "  -2.44%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").Value = "6.368"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").Value = "16.07"
$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("D23").Value = "11.72"
$ws.Range("E23").Value = "  -2.84%  "

$ws.Range("D24").Value = "22.172.70"
$ws.Range("E24").Value = "  -1.19%  "

$ws.Range("D25").Value = "2.291"
$ws.Range("E25").Value = "  -4.27%  "

$ws.Range("D26").Value = "2.534"
$ws.Range("E26").Value = "  -0.88%  "

$ws.Range("D27").Value = "150.65"
$ws.Range("E27").Value = "  -2.05%  "

$ws.Range("D28").Value = "19.10"
$ws.Range("E28").Value = "  -2.78%  "

$ws.Range("D29").Value = "4.933"
$ws.Range("E29").Value = "  -1.71%  "

$ws.Range("D30").Value = "1.757.03"
$ws.Range("E30").Value = "  +0.28%  "

$ws.Range("D31").Value = "122.19"
$ws.Range("E31").Value = "  -1.38%  "

$ws.Range("E32").Value = "  +3.42%  "

$ws.Range("D33").Value = "5.875"
$ws.Range("E33").Value = "  -4.05%  "

$ws.Range("D34").Value = "1.902"
$ws.Range("E34").Value = "  -5.59%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.08206"
$ws.Range("E35").Value = "  -1.13%  "

$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "9.326"
$ws.Range("E36").Value = "  -2.90%  "

$ws.Range("D37").Value = "0.06288"
$ws.Range("E37").Value = "  -1.40%  "

$ws.Range("D38").Value = "0.02323"
$ws.Range("E38").Value = "  -5.26%  "

$ws.Range("D39").Value = "5.297"
$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("D40").Value = "0.2147"
$ws.Range("E40").Value = "  -5.36%  "

$ws.Range("D41").Value = "1.233"
$ws.Range("E41").Value = "  -4.32%  "

$ws.Range("D42").Value = "10.97"
$ws.Range("E42").Value = "  -2.26%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.6056"
$ws.Range("E43").Value = "  -3.43%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Value = "13.78"
$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("D46").Value = "3.733"
$ws.Range("E46").Value = "  -1.39%  "

$ws.Range("D47").Value = "0.5843"
$ws.Range("E47").Value = "  -4.56%  "

$ws.Range("D48").Value = "1.967"
$ws.Range("E48").Value = "  -3.65%  "

$ws.Range("D49").Value = "121.55"
$ws.Range("E49").Value = "  -3.02%  "

$ws.Range("D50").Value = "1.172"
$ws.Range("E50").Value = "  -3.35%  "

$ws.Range("D51").Value = "0.07025"
$ws.Range("E51").Value = "  -2.94%  "
